$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 132
$ws.Range("B132").Value = 7483081
$ws.Range("C132").Value = "Ecuador LigaPro Serie A"
$ws.Range("D132").Value = "Ecuador LigaPro Serie A"
$ws.Range("E132").Value = 45255.83333333334
$ws.Range("F132").Value = "Deportivo Cuenca"
$ws.Range("G132").Value = "El Nacional"
$ws.Range("H132").Value = 1
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = "H"
$ws.Range("K132").Value = 2.75
$ws.Range("L132").Value = 3.25
$ws.Range("M132").Value = 2.55
$ws.Range("N132").Value = 3
$ws.Range("O132").Value = 3.3
$ws.Range("P132").Value = 2.3
$ws.Range("Q132").Value = 0.25
$ws.Range("R132").Value = 1.825
$ws.Range("S132").Value = 1.975
$ws.Range("T132").Value = 2.75
$ws.Range("U132").Value = 2
$ws.Range("V132").Value = 1.8
$ws.Range("W132").Value = 2
$ws.Range("X132").Value = -1
$ws.Range("Y132").Value = -1
$ws.Range("Z132").Value = 0.825
$ws.Range("AA132").Value = -1
$ws.Range("AB132").Value = -1
$ws.Range("AC132").Value = 0.8

# Row 133
$ws.Range("B133").Value = 7483189
$ws.Range("C133").Value = "Ecuador LigaPro Serie A"
$ws.Range("D133").Value = "Ecuador LigaPro Serie A"
$ws.Range("E133").Value = 45255.83333333334
$ws.Range("F133").Value = "Independiente del Valle"
$ws.Range("G133").Value = "Orense"
$ws.Range("H133").Value = 2
$ws.Range("I133").Value = 2
$ws.Range("J133").Value = "D"
$ws.Range("K133").Value = 1.4
$ws.Range("L133").Value = 4.75
$ws.Range("M133").Value = 7
$ws.Range("N133").Value = 1.4
$ws.Range("O133").Value = 4.5
$ws.Range("P133").Value = 8
$ws.Range("Q133").Value = -1.25
$ws.Range("R133").Value = 1.875
$ws.Range("S133").Value = 1.925
$ws.Range("T133").Value = 2.5
$ws.Range("U133").Value = 1.925
$ws.Range("V133").Value = 1.875
$ws.Range("W133").Value = -1
$ws.Range("X133").Value = 3.5
$ws.Range("Y133").Value = -1
$ws.Range("Z133").Value = -1
$ws.Range("AA133").Value = 0.925
$ws.Range("AB133").Value = 0.925
$ws.Range("AC133").Value = -1

# Row 134
$ws.Range("B134").Value = 7482832
$ws.Range("C134").Value = "Ecuador LigaPro Serie A"
$ws.Range("D134").Value = "Ecuador LigaPro Serie A"
$ws.Range("E134").Value = 45256.83333333334
$ws.Range("F134").Value = "Barcelona Guayaquil"
$ws.Range("G134").Value = "Guayaquil City"
$ws.Range("H134").Value = 2
$ws.Range("I134").Value = 1
$ws.Range("J134").Value = "H"
$ws.Range("K134").Value = 1.363
$ws.Range("L134").Value = 5
$ws.Range("M134").Value = 7.5
$ws.Range("N134").Value = 1.444
$ws.Range("O134").Value = 4
$ws.Range("P134").Value = 8
$ws.Range("Q134").Value = -1.25
$ws.Range("R134").Value = 2.05
$ws.Range("S134").Value = 1.75
$ws.Range("T134").Value = 2.5
$ws.Range("U134").Value = 1.95
$ws.Range("V134").Value = 1.85
$ws.Range("W134").Value = 0.444
$ws.Range("X134").Value = -1
$ws.Range("Y134").Value = -1
$ws.Range("Z134").Value = -0.5
$ws.Range("AA134").Value = 0.375
$ws.Range("AB134").Value = 0.95
$ws.Range("AC134").Value = -1

# Row 135
$ws.Range("B135").Value = 7483306
$ws.Range("C135").Value = "Ecuador LigaPro Serie A"
$ws.Range("D135").Value = "Ecuador LigaPro Serie A"
$ws.Range("E135").Value = 45256.83333333334
$ws.Range("F135").Value = "Tecnico Universitario"
$ws.Range("G135").Value = "Club Atletico Libertad"
$ws.Range("H135").Value = 1
$ws.Range("I135").Value = 1
$ws.Range("J135").Value = "D"
$ws.Range("K135").Value = 1.5
$ws.Range("L135").Value = 4.333
$ws.Range("M135").Value = 5.75
$ws.Range("N135").Value = 1.533
$ws.Range("O135").Value = 4.2
$ws.Range("P135").Value = 5.5
$ws.Range("Q135").Value = -1
$ws.Range("R135").Value = 1.925
$ws.Range("S135").Value = 1.875
$ws.Range("T135").Value = 2.25
$ws.Range("U135").Value = 1.8
$ws.Range("V135").Value = 2
$ws.Range("W135").Value = -1
$ws.Range("X135").Value = 3.2
$ws.Range("Y135").Value = -1
$ws.Range("Z135").Value = -1
$ws.Range("AA135").Value = 0.875
$ws.Range("AB135").Value = -0.5
$ws.Range("AC135").Value = 0.5

# Row 139
$ws.Range("B139").Value = 7528849
$ws.Range("C139").Value = "Ecuador LigaPro Serie A"
$ws.Range("D139").Value = "Ecuador LigaPro Serie A"
$ws.Range("E139").Value = 45262.70833333334
$ws.Range("F139").Value = "Guayaquil City"
$ws.Range("G139").Value = "Gualaceo SC"
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 2
$ws.Range("J139").Value = "A"
$ws.Range("K139").Value = 1.833
$ws.Range("L139").Value = 3.5
$ws.Range("M139").Value = 3.75
$ws.Range("N139").Value = 2.15
$ws.Range("O139").Value = 3.4
$ws.Range("P139").Value = 3
$ws.Range("Q139").Value = -0.25
$ws.Range("R139").Value = 1.825
$ws.Range("S139").Value = 1.975
$ws.Range("T139").Value = 2.5
$ws.Range("U139").Value = 1.85
$ws.Range("V139").Value = 1.95
$ws.Range("W139").Value = -1
$ws.Range("X139").Value = -1
$ws.Range("Y139").Value = 2
$ws.Range("Z139").Value = -1
$ws.Range("AA139").Value = 0.9750000000000001
$ws.Range("AB139").Value = -1
$ws.Range("AC139").Value = 0.95

# Row 140
$ws.Range("B140").Value = 7528859
$ws.Range("C140").Value = "Ecuador LigaPro Serie A"
$ws.Range("D140").Value = "Ecuador LigaPro Serie A"
$ws.Range("E140").Value = 45262.70833333334
$ws.Range("F140").Value = "Club Atletico Libertad"
$ws.Range("G140").Value = "Cumbaya FC"
$ws.Range("H140").Value = 3
$ws.Range("I140").Value = 1
$ws.Range("J140").Value = "H"
$ws.Range("K140").Value = 1.727
$ws.Range("L140").Value = 3.5
$ws.Range("M140").Value = 4.333
$ws.Range("N140").Value = 1.4
$ws.Range("O140").Value = 4.2
$ws.Range("P140").Value = 7
$ws.Range("Q140").Value = -1.25
$ws.Range("R140").Value = 2
$ws.Range("S140").Value = 1.8
$ws.Range("T140").Value = 2.5
$ws.Range("U140").Value = 1.95
$ws.Range("V140").Value = 1.85
$ws.Range("W140").Value = 0.3999999999999999
$ws.Range("X140").Value = -1
$ws.Range("Y140").Value = -1
$ws.Range("Z140").Value = 1
$ws.Range("AA140").Value = -1
$ws.Range("AB140").Value = 0.95
$ws.Range("AC140").Value = -1

# Row 142
$ws.Range("B142").Value = 7528848
$ws.Range("C142").Value = "Ecuador LigaPro Serie A"
$ws.Range("D142").Value = "Ecuador LigaPro Serie A"
$ws.Range("E142").Value = 45263.83333333334
$ws.Range("F142").Value = "Emelec"
$ws.Range("G142").Value = "Deportivo Cuenca"
$ws.Range("H142").Value = 2
$ws.Range("I142").Value = 1
$ws.Range("J142").Value = "H"
$ws.Range("K142").Value = 1.75
$ws.Range("L142").Value = 3.5
$ws.Range("M142").Value = 4.2
$ws.Range("N142").Value = 2.4
$ws.Range("O142").Value = 3.1
$ws.Range("P142").Value = 2.75
$ws.Range("Q142").Value = -0.25
$ws.Range("R142").Value = 2.05
$ws.Range("S142").Value = 1.75
$ws.Range("T142").Value = 2.25
$ws.Range("U142").Value = 1.8
$ws.Range("V142").Value = 2
$ws.Range("W142").Value = 1.4
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = -1
$ws.Range("Z142").Value = 1.05
$ws.Range("AA142").Value = -1
$ws.Range("AB142").Value = 0.8
$ws.Range("AC142").Value = -1

# Row 144
$ws.Range("B144").Value = 7528852
$ws.Range("C144").Value = "Ecuador LigaPro Serie A"
$ws.Range("D144").Value = "Ecuador LigaPro Serie A"
$ws.Range("E144").Value = 45263.83333333334
$ws.Range("F144").Value = "Delfin SC"
$ws.Range("G144").Value = "Tecnico Universitario"
$ws.Range("H144").Value = 2
$ws.Range("I144").Value = 2
$ws.Range("J144").Value = "D"
$ws.Range("K144").Value = 2.1
$ws.Range("L144").Value = 3.4
$ws.Range("M144").Value = 3.1
$ws.Range("N144").Value = 2.1
$ws.Range("O144").Value = 3.4
$ws.Range("P144").Value = 3.1
$ws.Range("Q144").Value = -0.25
$ws.Range("R144").Value = 1.8
$ws.Range("S144").Value = 2
$ws.Range("T144").Value = 2.25
$ws.Range("U144").Value = 1.9
$ws.Range("V144").Value = 1.9
$ws.Range("W144").Value = -1
$ws.Range("X144").Value = 2.4
$ws.Range("Y144").Value = -1
$ws.Range("Z144").Value = -0.5
$ws.Range("AA144").Value = 0.5
$ws.Range("AB144").Value = 0.8999999999999999
$ws.Range("AC144").Value = -1

# Row 145
$ws.Range("B145").Value = 7528858
$ws.Range("C145").Value = "Ecuador LigaPro Serie A"
$ws.Range("D145").Value = "Ecuador LigaPro Serie A"
$ws.Range("E145").Value = 45263.83333333334
$ws.Range("F145").Value = "Orense"
$ws.Range("G145").Value = "SD Aucas"
$ws.Range("H145").Value = 1
$ws.Range("I145").Value = 2
$ws.Range("J145").Value = "A"
$ws.Range("K145").Value = 2.2
$ws.Range("L145").Value = 3.2
$ws.Range("M145").Value = 3.2
$ws.Range("N145").Value = 1.95
$ws.Range("O145").Value = 3.2
$ws.Range("P145").Value = 3.8
$ws.Range("Q145").Value = -0.5
$ws.Range("R145").Value = 1.95
$ws.Range("S145").Value = 1.85
$ws.Range("T145").Value = 2.25
$ws.Range("U145").Value = 1.85
$ws.Range("V145").Value = 1.95
$ws.Range("W145").Value = -1
$ws.Range("X145").Value = -1
$ws.Range("Y145").Value = 2.8
$ws.Range("Z145").Value = -1
$ws.Range("AA145").Value = 0.8500000000000001
$ws.Range("AB145").Value = 0.8500000000000001
$ws.Range("AC145").Value = -1

# Row 200
$ws.Range("B200").Value = 8069719
$ws.Range("C200").Value = "Ecuador LigaPro Serie A"
$ws.Range("D200").Value = "Ecuador LigaPro Serie A"
$ws.Range("E200").Value = 45395.72916666666
$ws.Range("F200").Value = "Macara"
$ws.Range("G200").Value = "Orense"
$ws.Range("K200").Value = 1.95
$ws.Range("L200").Value = 3.25
$ws.Range("M200").Value = 3.5
$ws.Range("N200").Value = 1.571
$ws.Range("O200").Value = 3.6
$ws.Range("P200").Value = 5.25
$ws.Range("Q200").Value = -0.75
$ws.Range("R200").Value = 1.775
$ws.Range("S200").Value = 2.025
$ws.Range("T200").Value = 2.5
$ws.Range("U200").Value = 1.975
$ws.Range("V200").Value = 1.825
$ws.Range("W200").Value = 0
$ws.Range("X200").Value = 0
$ws.Range("Y200").Value = 0
$ws.Range("Z200").Value = 0
$ws.Range("AA200").Value = 0

# Row 201
$ws.Range("B201").Value = 8069537
$ws.Range("C201").Value = "Ecuador LigaPro Serie A"
$ws.Range("D201").Value = "Ecuador LigaPro Serie A"
$ws.Range("E201").Value = 45395.83333333334
$ws.Range("F201").Value = "Emelec"
$ws.Range("G201").Value = "Cumbaya FC"
$ws.Range("K201").Value = 1.28
$ws.Range("L201").Value = 5.5
$ws.Range("M201").Value = 8.5
$ws.Range("N201").Value = 1.25
$ws.Range("O201").Value = 5.5
$ws.Range("P201").Value = 10
$ws.Range("Q201").Value = -1.75
$ws.Range("R201").Value = 1.95
$ws.Range("S201").Value = 1.85
$ws.Range("T201").Value = 2.75
$ws.Range("U201").Value = 1.95
$ws.Range("V201").Value = 1.85
$ws.Range("W201").Value = 0
$ws.Range("X201").Value = 0
$ws.Range("Y201").Value = 0
$ws.Range("Z201").Value = 0
$ws.Range("AA201").Value = 0

# Row 202
$ws.Range("B202").Value = 7773067
$ws.Range("C202").Value = "Ecuador LigaPro Serie A"
$ws.Range("D202").Value = "Ecuador LigaPro Serie A"
$ws.Range("E202").Value = 45396.625
$ws.Range("F202").Value = "Universidad Catolica del Ecuador"
$ws.Range("G202").Value = "Independiente del Valle"
$ws.Range("K202").Value = 2.3
$ws.Range("L202").Value = 3.1
$ws.Range("M202").Value = 3
$ws.Range("N202").Value = 2.6
$ws.Range("O202").Value = 3.1
$ws.Range("P202").Value = 2.6
$ws.Range("Q202").Value = 0
$ws.Range("R202").Value = 1.875
$ws.Range("S202").Value = 1.925
$ws.Range("T202").Value = 2.25
$ws.Range("U202").Value = 1.825
$ws.Range("V202").Value = 1.975
$ws.Range("W202").Value = 0
$ws.Range("X202").Value = 0
$ws.Range("Y202").Value = 0
$ws.Range("Z202").Value = 0
$ws.Range("AA202").Value = 0

# Row 203
$ws.Range("B203").Value = 8069720
$ws.Range("C203").Value = "Ecuador LigaPro Serie A"
$ws.Range("D203").Value = "Ecuador LigaPro Serie A"
$ws.Range("E203").Value = 45396.72916666666
$ws.Range("F203").Value = "El Nacional"
$ws.Range("G203").Value = "Imbabura"
$ws.Range("K203").Value = 1.7
$ws.Range("L203").Value = 3.6
$ws.Range("M203").Value = 4.2
$ws.Range("N203").Value = 1.615
$ws.Range("O203").Value = 3.75
$ws.Range("P203").Value = 4.5
$ws.Range("Q203").Value = -0.75
$ws.Range("R203").Value = 1.825
$ws.Range("S203").Value = 1.975
$ws.Range("T203").Value = 2.5
$ws.Range("U203").Value = 1.95
$ws.Range("V203").Value = 1.85
$ws.Range("W203").Value = 0
$ws.Range("X203").Value = 0
$ws.Range("Y203").Value = 0
$ws.Range("Z203").Value = 0
$ws.Range("AA203").Value = 0

# Row 204
$ws.Range("B204").Value = 8069721
$ws.Range("C204").Value = "Ecuador LigaPro Serie A"
$ws.Range("D204").Value = "Ecuador LigaPro Serie A"
$ws.Range("E204").Value = 45396.83333333334
$ws.Range("F204").Value = "Club Atletico Libertad"
$ws.Range("G204").Value = "Barcelona Guayaquil"
$ws.Range("K204").Value = 4
$ws.Range("L204").Value = 3.4
$ws.Range("M204").Value = 1.85
$ws.Range("N204").Value = 4.2
$ws.Range("O204").Value = 3.4
$ws.Range("P204").Value = 1.8
$ws.Range("Q204").Value = 0.5
$ws.Range("R204").Value = 2.025
$ws.Range("S204").Value = 1.775
$ws.Range("T204").Value = 2.5
$ws.Range("U204").Value = 1.975
$ws.Range("V204").Value = 1.825
$ws.Range("W204").Value = 0
$ws.Range("X204").Value = 0
$ws.Range("Y204").Value = 0
$ws.Range("Z204").Value = 0
$ws.Range("AA204").Value = 0

# Row 205
$ws.Range("B205").Value = 7773504
$ws.Range("C205").Value = "Ecuador LigaPro Serie A"
$ws.Range("D205").Value = "Ecuador LigaPro Serie A"
$ws.Range("E205").Value = 45397.875
$ws.Range("F205").Value = "LDU Quito"
$ws.Range("G205").Value = "Delfin SC"
$ws.Range("K205").Value = 1.363
$ws.Range("L205").Value = 4.333
$ws.Range("M205").Value = 7.5
$ws.Range("N205").Value = 1.333
$ws.Range("O205").Value = 4.5
$ws.Range("P205").Value = 8
$ws.Range("Q205").Value = -1.5
$ws.Range("R205").Value = 2
$ws.Range("S205").Value = 1.8
$ws.Range("T205").Value = 2.75
$ws.Range("U205").Value = 1.9
$ws.Range("V205").Value = 1.9
$ws.Range("W205").Value = 0
$ws.Range("X205").Value = 0
$ws.Range("Y205").Value = 0
$ws.Range("Z205").Value = 0
$ws.Range("AA205").Value = 0

# Remove the now-obsolete last row (match 7773504 has been folded into row 205's
# slot as part of the refreshed fixture list; the trailing duplicate row is deleted).
$ws.Rows.Item(206).Delete()
